# Daily attendance processing - 2025-10-23 03:39:36
# Re-sorts the comma-separated "Recorded By" (column G) values on each
# attendance row into alphabetical (case-insensitive) order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $original = $cell.Text

    if ($original -ne $null -and $original -ne "") {
        $parts = $original.Split(",")

        if ($parts.Length -gt 1) {
            $trimmedParts = @()
            foreach ($part in $parts) {
                $trimmedParts += $part.Trim()
            }

            $sortedParts = $trimmedParts | Sort-Object

            $updated = $sortedParts -join ", "

            if ($updated -ne $original) {
                $cell.Value = $updated
            }
        }
    }
}
